# Update the "想去人数" (number of people wanting to go) column (F) values
# across the "展览", "本地生活" and "全部类型" worksheets, per the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 641
    4  = 41
    5  = 1969
    6  = 5624
    7  = 1580
    9  = 3197
    11 = 44
    12 = 1331
    13 = 4461
    14 = 1064
    15 = 1692
    18 = 43
    19 = 45
    20 = 163
    21 = 154
    22 = 1008
    23 = 302
    28 = 4
    30 = 399
    31 = 76
    33 = 347
    34 = 175
    36 = 1722
    38 = 1044
    40 = 265
    41 = 624
    42 = 347
    43 = 24
    45 = 23
    46 = 425
    47 = 375
    49 = 145
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 774

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 774
    4  = 41
    5  = 5624
    6  = 1580
    9  = 3197
    10 = 1331
    11 = 4461
    12 = 1064
    16 = 43
    19 = 45
    20 = 163
    21 = 154
    23 = 1008
    24 = 302
    29 = 4
    31 = 399
    33 = 175
    35 = 1722
    37 = 1044
    41 = 265
    42 = 347
    44 = 425
    45 = 375
    48 = 145
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
